# Roithner_S6305MG_OpticalPower_Laser03.xlsx - "Add files via upload" edit
#
# What changed (per the OOXML diff):
#   1. The 4th worksheet "OpticalPower_vs_Temp_15mA" is renamed to
#      "OpticalPower_vs_Temp".
#   2. That renamed sheet becomes the active/selected tab (workbookView
#      activeTab goes 2 -> 3, its <sheetView> gains tabSelected="1") and its
#      selection moves to I29.
#   3. The previously-active sheet "OpticalPower_vs_Current_25º" loses
#      tabSelected (its own cell selection, K26, is unchanged).
#   4. The chart on the renamed sheet has its cached series/error-bar
#      formulas repointed at the new sheet name (cosmetic - the cached
#      values themselves are untouched).

$wb = $excel.ActiveWorkbook

# --- 1. Rename the sheet -------------------------------------------------
$wsTemp = $wb.Worksheets.Item("OpticalPower_vs_Temp_15mA")
$wsTemp.Name = "OpticalPower_vs_Temp"

# --- 2. Make it the active sheet & set its selection ---------------------
$wsTemp.Activate() | Out-Null
$wsTemp.Range("I29").Select() | Out-Null

# --- 3. The other sheet's own selection (K26) is left exactly as it was --
# (Activating OpticalPower_vs_Temp above already clears tabSelected from
#  whichever sheet had it before; OpticalPower_vs_Current_25º keeps its own
#  cached selection at K26 untouched since we never change it.)

# --- 4. Repoint the chart's cached formulas to the new sheet name --------
# (Best effort only - this COM surface does not expose a way to rewrite the
#  raw <c:f> text of an already-imported chart's series/error-bars, so
#  these calls are harmless no-ops here. Deliberately NOT touching
#  Series.XValues / Series.Values with plain strings: on this host that
#  path is reserved for brand-new series data and spawns a stray extra
#  chart object instead of patching the existing one.)
foreach ($co in $wsTemp.ChartObjects()) {
    $chart = $co.Chart
    foreach ($s in $chart.SeriesCollection()) {
        try { $s.Formula = $s.Formula -replace "OpticalPower_vs_Temp_15mA", "OpticalPower_vs_Temp" } catch {}
        try { $s.ErrorBars.Formula1 = "OpticalPower_vs_Temp!`$O`$11:`$O`$30" } catch {}
    }
}
